$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / value changes -------------------------------------------------
$ws.Range("B2").Value = "rohan"
$ws.Range("J2").Value = "tintu@gmail.com"
$ws.Range("L2").Value = "Reading ,Drawing"

$ws.Range("B3").Value = "mini"
$ws.Range("H3").Value = "abcd"
$ws.Range("L3").Value = "Reading ,Writing"

# --- Font colour change (theme colour -> explicit black) ------------------
$ws.Range("I2:I3").Font.Color = 0
$ws.Range("K2:K3").Font.Color = 0

# --- Row height changes -----------------------------------------------------
$ws.Rows.Item(1).RowHeight = 19.5
$ws.Rows.Item(2).RowHeight = 19.5
$ws.Rows.Item(3).RowHeight = 19.5
